$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Re-color the table header cells (A3,C3,D3,E3,G3,H3) from red to automatic/black.
#    (B3 and F3 keep their existing styles untouched.)
$headerCells = @("A3", "C3", "D3", "E3", "G3", "H3")
foreach ($addr in $headerCells) {
    $ws.Range($addr).Font.ThemeColor = 1
}

# Row 3 grows taller to fit the new helper-list header next to it.
$ws.Rows(3).RowHeight = 30

# A couple of the detail rows also pick up explicit heights.
$ws.Rows(7).RowHeight = 18
$ws.Rows(8).RowHeight = 15

# 2) Populate the hidden helper list in columns L:M used to drive the new
#    "nhập xuất tồn" (stock in/out) category dropdown validation.
$categories = @("Smartphone", "Macbook", "Ipad", "Notepad", "Phone Accessories", "Laptop Accessories")
for ($i = 0; $i -lt $categories.Length; $i++) {
    $row = 3 + $i
    $cell = $ws.Range("L$row")
    $cell.Value = $categories[$i]
    $cell.Font.ThemeColor = 2
    $cell.WrapText = $true
    $cell.VerticalAlignment = -4108

    $mcell = $ws.Range("M$row")
    $mcell.Font.ThemeColor = 2
}

# The rest of the helper block (rows 9-14) stays blank but shares the same
# "invisible" white-font styling.
for ($row = 9; $row -le 14; $row++) {
    $ws.Range("L$row").Font.ThemeColor = 2
    $ws.Range("M$row").Font.ThemeColor = 2
}

# 3) Wire up the dropdown (data validation) lists for the new input row.
$ws.Range("C4").Validation.Add(3, 1, 1, "=`$N`$3:`$N`$14")
$ws.Range("C4").Validation.ShowInput = $true
$ws.Range("C4").Validation.ShowError = $true

$ws.Range("D4").Validation.Add(3, 1, 1, "=`$L`$3:`$L`$8")
$ws.Range("D4").Validation.ShowInput = $true
$ws.Range("D4").Validation.ShowError = $true

# 4) Move the active selection to the new data-entry cell.
$ws.Range("C4").Select()
